# Refresh the "cryptos" price/volume table with the latest scraped values.
# Price (col D) and Volume(1h) (col E) are stored as plain text, so any
# value that could be mistaken for a number gets a Text number format
# before it is written and the cell style is reset back to "Normal"
# afterwards (this keeps the cell as text without leaving a stray
# quote-prefixed / Text-formatted style behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.391.92"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).Value = "  -1.24%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.532.65"
$ws.Cells.Item(3, 4).Style = "Normal"

$ws.Cells.Item(3, 5).Value = "  -2.60%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "  +0.08%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "196.54"
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(5, 5).Value = "  +0.14%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "582.31"
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(6, 5).Value = "  -3.65%  "

$ws.Cells.Item(7, 5).Value = "  -2.76%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.625"
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(10, 5).Value = "  -4.03%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "51.73"
$ws.Cells.Item(11, 4).Style = "Normal"

$ws.Cells.Item(11, 5).Value = "  -4.36%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000286"
$ws.Cells.Item(12, 4).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "  -6.81%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "9.26"
$ws.Cells.Item(13, 4).Style = "Normal"

$ws.Cells.Item(13, 5).Value = "  -3.16%  "

$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 2).Style = "Normal"

$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(14, 3).Style = "Normal"

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.087.77"
$ws.Cells.Item(14, 4).Style = "Normal"

$ws.Cells.Item(14, 5).Value = "  -2.70%  "

$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "BitcoinCash"
$ws.Cells.Item(15, 2).Style = "Normal"

$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(15, 3).Style = "Normal"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "663.23"
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(15, 5).Value = "  +11.48%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "69.443.07"
$ws.Cells.Item(16, 4).Style = "Normal"

$ws.Cells.Item(16, 5).Value = "  -1.43%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.528.73"
$ws.Cells.Item(17, 4).Style = "Normal"

$ws.Cells.Item(17, 5).Value = "  -3.10%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "12.43"
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(18, 5).Value = "  -6.41%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "18.53"
$ws.Cells.Item(19, 4).Style = "Normal"

$ws.Cells.Item(19, 5).Value = "  -3.66%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.968"
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(21, 5).Value = "  -3.12%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "18.40"
$ws.Cells.Item(22, 4).Style = "Normal"

$ws.Cells.Item(22, 5).Value = "  +3.85%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "105.29"
$ws.Cells.Item(23, 4).Style = "Normal"

$ws.Cells.Item(23, 5).Value = "  +2.82%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.26"
$ws.Cells.Item(24, 4).Style = "Normal"

$ws.Cells.Item(24, 5).Value = "  +2.18%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "4.38"
$ws.Cells.Item(25, 4).Style = "Normal"

$ws.Cells.Item(25, 5).Value = "  -5.50%  "

$ws.Cells.Item(26, 5).Value = "  -3.76%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.19"
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(27, 5).Value = "  -5.97%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.68"
$ws.Cells.Item(28, 4).Style = "Normal"

$ws.Cells.Item(28, 5).Value = "  +0.43%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "33.36"
$ws.Cells.Item(29, 4).Style = "Normal"

$ws.Cells.Item(29, 5).Value = "  -2.18%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.41"
$ws.Cells.Item(30, 4).Style = "Normal"

$ws.Cells.Item(30, 5).Value = "  -7.45%  "

$ws.Cells.Item(31, 5).Value = "  -5.04%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "11.88"
$ws.Cells.Item(32, 4).Style = "Normal"

$ws.Cells.Item(32, 5).Value = "  -3.94%  "

$ws.Cells.Item(33, 5).Value = "  -5.63%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "61.96"
$ws.Cells.Item(34, 4).Style = "Normal"

$ws.Cells.Item(34, 5).Value = "  -2.12%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.794.84"
$ws.Cells.Item(35, 4).Style = "Normal"

$ws.Cells.Item(35, 5).Value = "  -3.26%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0₃0812"
$ws.Cells.Item(36, 4).Style = "Normal"

$ws.Cells.Item(36, 5).Value = "  -10.93%  "

$ws.Cells.Item(37, 5).Value = "  +0.07%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.67"
$ws.Cells.Item(38, 4).Style = "Normal"

$ws.Cells.Item(38, 5).Value = "  +3.30%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "504.25"
$ws.Cells.Item(39, 4).Style = "Normal"

$ws.Cells.Item(39, 5).Value = "  -4.82%  "

$ws.Cells.Item(40, 5).Value = "  -7.04%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.372"
$ws.Cells.Item(41, 4).Style = "Normal"

$ws.Cells.Item(41, 5).Value = "  -5.43%  "

$ws.Cells.Item(42, 5).Value = "  -0.16%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "34.57"
$ws.Cells.Item(43, 4).Style = "Normal"

$ws.Cells.Item(43, 5).Value = "  -7.99%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0456"
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(44, 5).Value = "  -0.14%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.89"
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(45, 5).Value = "  +0.67%  "

$ws.Cells.Item(46, 5).Value = "  +0.43%  "

$ws.Cells.Item(47, 5).Value = "  -3.32%  "

$ws.Cells.Item(48, 5).Value = "  -0.32%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.30"
$ws.Cells.Item(49, 4).Style = "Normal"

$ws.Cells.Item(49, 5).Value = "  -3.96%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.77"
$ws.Cells.Item(50, 4).Style = "Normal"

$ws.Cells.Item(50, 5).Value = "  +69.15%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.77"
$ws.Cells.Item(51, 4).Style = "Normal"

$ws.Cells.Item(51, 5).Value = "  +18.53%  "
